$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.087.10'
$ws.Range('E2').Value = '  +0.79%  '

$ws.Range('D3').Value = '2.637.34'
$ws.Range('E3').Value = '  +3.01%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').Value = "'523.03"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.66%  '

$ws.Range('D6').Value = "'146.59"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.41%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = "'0.572"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.38%  '

$ws.Range('D9').Value = '2.652.37'
$ws.Range('E9').Value = '  +3.05%  '

$ws.Range('D10').Value = "'6.33"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.06%  '

$ws.Range('E11').Value = '  +1.94%  '

$ws.Range('E12').Value = '  +0.99%  '

$ws.Range('E13').Value = '  -1.13%  '

$ws.Range('D14').Value = '3.099.87'
$ws.Range('E14').Value = '  +2.99%  '

$ws.Range('D15').Value = '59.066.79'
$ws.Range('E15').Value = '  +0.82%  '

$ws.Range('D16').Value = "'21.05"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.45%  '

$ws.Range('E17').Value = '  +0.64%  '

$ws.Range('D18').Value = '2.636.53'
$ws.Range('E18').Value = '  +2.86%  '

$ws.Range('D19').Value = "'347.32"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.24%  '

$ws.Range('E20').Value = '  -0.65%  '

$ws.Range('D21').Value = "'10.30"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.77%  '

$ws.Range('E22').Value = '  +2.78%  '

$ws.Range('E23').Value = '  -0.27%  '

$ws.Range('D24').Value = "'61.87"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.02%  '

$ws.Range('E25').Value = '  +1.27%  '

$ws.Range('D26').Value = "'0.166"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.93%  '

$ws.Range('E27').Value = '  -0.13%  '

$ws.Range('D28').Value = '0.0₃0808'
$ws.Range('E28').Value = '  +1.03%  '

$ws.Range('E29').Value = '  +1.80%  '

$ws.Range('E30').Value = '  -0.04%  '

$ws.Range('D31').Value = "'6.28"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +5.28%  '

$ws.Range('E32').Value = '  +3.07%  '

$ws.Range('D33').Value = "'18.98"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.87%  '

$ws.Range('D34').Value = "'149.96"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.32%  '

$ws.Range('D35').Value = "'0.983"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +9.16%  '

$ws.Range('D36').Value = "'4.01"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.95%  '

$ws.Range('E37').Value = '  +1.58%  '

$ws.Range('D38').Value = "'36.77"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.10%  '

$ws.Range('E39').Value = '  +0.80%  '

$ws.Range('E40').Value = '  +3.57%  '

$ws.Range('E41').Value = '  +1.57%  '

$ws.Range('D42').Value = "'279.34"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.74%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'0.995"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.09%  '

$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = "'0.610"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.92%  '

$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').Value = "'0.0986"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.06%  '

$ws.Range('E46').Value = '  +3.09%  '

$ws.Range('E47').Value = '  -1.88%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = "'0.0230"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.53%  '

$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').Value = "'10.30"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.43%  '

$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '1.990.59'
$ws.Range('E50').Value = '  +3.69%  '

$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = "'4.67"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.21%  '
